$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.568.31"
$ws.Range("E2").Value = "  +4.13%  "
$ws.Range("D3").Value = "2.331.88"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'547.10"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").Value = "'131.54"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").Value = "2.329.93"
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").Value = "'5.51"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "'0.337"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").Value = "'23.84"
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("D15").Value = "2.747.86"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").Value = "60.496.33"
$ws.Range("E16").Value = "  +4.16%  "
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "2.336.65"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").Value = "'10.65"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").Value = "'315.17"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "'6.65"
$ws.Range("E22").Value = "  +2.59%  "
$ws.Range("D23").Value = "'0.996"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "'64.17"
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("D25").Value = "'0.170"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'7.90"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").Value = "'1.35"
$ws.Range("E28").Value = "  +4.62%  "
$ws.Range("D29").Value = "'1.19"
$ws.Range("E29").Value = "  +11.48%  "
$ws.Range("D30").Value = "'173.01"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("E31").Value = "  +2.68%  "
$ws.Range("D32").Value = "0.0₃0736"
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("D33").Value = "'5.96"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("E34").Value = "  +11.79%  "
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").Value = "'17.91"
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'4.09"
$ws.Range("E39").Value = "  +4.56%  "
$ws.Range("D40").Value = "'329.31"
$ws.Range("E40").Value = "  +14.18%  "
$ws.Range("D41").Value = "'1.54"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("D42").Value = "'38.00"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").Value = "'139.12"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("D45").Value = "'0.0945"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").Value = "'19.40"
$ws.Range("E46").Value = "  +7.18%  "
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("D48").Value = "'0.562"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0214"
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0220"
$ws.Range("E50").Value = "  +19.97%  "
$ws.Range("D51").Value = "'11.02"
$ws.Range("E51").Value = "  +0.68%  "
